# Refresh the NATMI ligand/receptor (Fcer2a-Cr2) TPM-derived metrics on the
# active worksheet. The underlying per-cluster TPM values were recomputed
# ("update scripts wuth new tpm"), which changes the ligand/receptor average &
# total expression values (columns G,H,M,N), their derived-specificity
# fractions (columns I,J,O,P), and the edge weights/specificities that are
# multiplicatively derived from them (columns Q,R,S,T). Columns A-F,K,L stay
# untouched since the expressing-cell counts/detection rates didn't change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.260875
$ws.Range("H2").Value = 0.7826249999999999
$ws.Range("I2").Value = 0.2304058263197128
$ws.Range("J2").Value = 0.2304058263197127
$ws.Range("M2").Value = 0.659041
$ws.Range("N2").Value = 1.977123
$ws.Range("O2").Value = 0.3819216825298216
$ws.Range("P2").Value = 0.3819216825298216
$ws.Range("Q2").Value = 0.171927320875
$ws.Range("R2").Value = 1.547345887875
$ws.Range("S2").Value = 0.08799698085269855
$ws.Range("T2").Value = 0.08799698085269854
# Row 3
$ws.Range("G3").Value = 0.260875
$ws.Range("H3").Value = 0.7826249999999999
$ws.Range("I3").Value = 0.2304058263197128
$ws.Range("J3").Value = 0.2304058263197127
$ws.Range("O3").Value = 0.3628134576423628
$ws.Range("P3").Value = 0.3628134576423628
$ws.Range("Q3").Value = 0.1633254895
$ws.Range("R3").Value = 1.4699294055
$ws.Range("S3").Value = 0.0835943345080007
$ws.Range("T3").Value = 0.08359433450800069
# Row 4
$ws.Range("G4").Value = 0.260875
$ws.Range("H4").Value = 0.7826249999999999
$ws.Range("I4").Value = 0.2304058263197128
$ws.Range("J4").Value = 0.2304058263197127
$ws.Range("O4").Value = 0.2552648598278156
$ws.Range("P4").Value = 0.2552648598278155
$ws.Range("Q4").Value = 0.114911002625
$ws.Range("R4").Value = 1.034199023625
$ws.Range("S4").Value = 0.05881451095901351
$ws.Range("T4").Value = 0.05881451095901349
# Row 5
$ws.Range("I5").Value = 0.5522653592108161
$ws.Range("J5").Value = 0.5522653592108161
$ws.Range("M5").Value = 0.659041
$ws.Range("N5").Value = 1.977123
$ws.Range("O5").Value = 0.3819216825298216
$ws.Range("P5").Value = 0.3819216825298216
$ws.Range("Q5").Value = 0.4120967995376666
$ws.Range("R5").Value = 3.708871195839
$ws.Range("S5").Value = 0.2109221151927312
$ws.Range("T5").Value = 0.2109221151927312
# Row 6
$ws.Range("I6").Value = 0.5522653592108161
$ws.Range("J6").Value = 0.5522653592108161
$ws.Range("O6").Value = 0.3628134576423628
$ws.Range("P6").Value = 0.3628134576423628
$ws.Range("S6").Value = 0.2003693045113777
$ws.Range("T6").Value = 0.2003693045113777
# Row 7
$ws.Range("I7").Value = 0.5522653592108161
$ws.Range("J7").Value = 0.5522653592108161
$ws.Range("O7").Value = 0.2552648598278156
$ws.Range("P7").Value = 0.2552648598278155
$ws.Range("S7").Value = 0.1409739395067072
$ws.Range("T7").Value = 0.1409739395067072
# Row 8
$ws.Range("G8").Value = 0.2460686666666667
$ws.Range("H8").Value = 0.738206
$ws.Range("I8").Value = 0.2173288144694712
$ws.Range("J8").Value = 0.2173288144694712
$ws.Range("M8").Value = 0.659041
$ws.Range("N8").Value = 1.977123
$ws.Range("O8").Value = 0.3819216825298216
$ws.Range("P8").Value = 0.3819216825298216
$ws.Range("Q8").Value = 0.1621693401486667
$ws.Range("R8").Value = 1.459524061338
$ws.Range("S8").Value = 0.08300258648439189
$ws.Range("T8").Value = 0.08300258648439188
# Row 9
$ws.Range("G9").Value = 0.2460686666666667
$ws.Range("H9").Value = 0.738206
$ws.Range("I9").Value = 0.2173288144694712
$ws.Range("J9").Value = 0.2173288144694712
$ws.Range("O9").Value = 0.3628134576423628
$ws.Range("P9").Value = 0.3628134576423628
$ws.Range("S9").Value = 0.07884981862298442
$ws.Range("T9").Value = 0.0788498186229844
# Row 10
$ws.Range("G10").Value = 0.2460686666666667
$ws.Range("H10").Value = 0.738206
$ws.Range("I10").Value = 0.2173288144694712
$ws.Range("J10").Value = 0.2173288144694712
$ws.Range("O10").Value = 0.2552648598278156
$ws.Range("P10").Value = 0.2552648598278155
$ws.Range("R10").Value = 0.9755015804939999
$ws.Range("S10").Value = 0.05547640936209491
$ws.Range("T10").Value = 0.05547640936209489

Write-Output "Applied 83 cell updates"